# Update doctor_MA average column (AF) values for rows 4-13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "AF4"  = 0.639
    "AF5"  = 0.833
    "AF6"  = 0.723
    "AF7"  = 0.785
    "AF8"  = 0.772
    "AF9"  = 0.667
    "AF10" = 0.833
    "AF11" = 0.833
    "AF12" = 1.2
    "AF13" = 1.667
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
